# Actualización desde MV -datos-
# Update existing rows 127, 137, 138 with new figures, and append daily
# subscription data for August 2021 (rows 147-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing cells -------------------------------------------------
$ws.Cells.Item(127, 2).Value  = 396   # B127
$ws.Cells.Item(127, 10).Value = 58    # J127

$ws.Cells.Item(137, 2).Value  = 298   # B137
$ws.Cells.Item(137, 9).Value  = 48    # I137

$ws.Cells.Item(138, 2).Value  = 399   # B138
$ws.Cells.Item(138, 9).Value  = 52    # I138

# --- 2. Append new daily rows (147-168) for August 2021 ----------------------
$dates = @(
    "02-08-2021","03-08-2021","04-08-2021","05-08-2021","06-08-2021",
    "09-08-2021","10-08-2021","11-08-2021","12-08-2021","13-08-2021",
    "16-08-2021","17-08-2021","18-08-2021","19-08-2021","20-08-2021",
    "23-08-2021","24-08-2021","25-08-2021","26-08-2021","27-08-2021",
    "30-08-2021","31-08-2021"
)

$values = @(
    @(48,0,0,0,0,16,0,0,32),
    @(82,0,0,0,39,0,19,8,16),
    @(454,0,321,0,78,0,20,21,14),
    @(346,0,186,0,77,27,30,4,22),
    @(184,0,0,0,24,77,29,33,21),
    @(373,0,0,156,16,0,146,53,2),
    @(455,96,0,84,166,26,56,0,28),
    @(633,0,246,105,38,57,107,37,43),
    @(359,0,0,0,115,90,155,0,0),
    @(235,0,0,0,40,53,63,39,40),
    @(695,0,92,0,0,314,244,43,3),
    @(221,96,0,0,0,16,76,32,0),
    @(513,191,0,0,40,16,214,14,38),
    @(183,0,0,32,24,0,68,15,46),
    @(136,0,0,0,0,0,60,73,3),
    @(302,96,0,0,0,22,161,9,15),
    @(412,0,249,63,24,0,71,0,6),
    @(630,153,0,156,200,0,86,34,2),
    @(779,31,156,146,366,20,22,32,6),
    @(649,158,81,31,368,0,0,10,0),
    @(377,0,48,0,250,0,79,0,0),
    @(320,0,205,31,0,0,36,45,3)
)

$startRow = 147

# The date labels must stay as plain text (matching the existing "DD-MM-YYYY"
# strings used throughout column A) instead of being auto-converted into
# date serial numbers, so force the number format to Text for that range
# before writing the values, then restore the style so the cells keep the
# workbook's default (unstyled) look.
$dateRange = $ws.Range("A$startRow`:A$($startRow + $dates.Length - 1)")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]

    $rowValues = $values[$i]
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($row, 2 + $col).Value = $rowValues[$col]
    }
}

$dateRange.Style = "Normal"
